# AGA206 Assessment 2 Checklist - apply "hazard, powerups and tilt" update
# - Clear a few stale notes (G12, G13, G18, G22, G23, G25, G29)
# - Update several task notes (G19, G20, G21, G24, G26, G27, G32)
# - Mark newly-completed optional tasks as Done by checking their linked
#   checkboxes (J20, J21, J26, J27, J32) which drives F (status) and K (points)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear stale / no-longer-relevant notes ---
$ws.Range("G12").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("G18").ClearContents()
$ws.Range("G22").ClearContents()
$ws.Range("G23").ClearContents()
$ws.Range("G25").ClearContents()
$ws.Range("G29").ClearContents()

# --- Update notes text ---
$ws.Range("G19").Value = "script done. "
$ws.Range("G20").Value = "all the assets complete"
$ws.Range("G21").Value = "one tilt level"
$ws.Range("G24").Value = "they move"
$ws.Range("G26").Value = "done"
$ws.Range("G27").Value = "this waas annoying"
$ws.Range("G32").Value = "three hazards"

# --- Check the checkboxes for the newly completed optional tasks ---
# (World Tilt Mode, New Look and Feel, Powerups -Speed/-Size Related, Moving Hazard)
$ws.Range("J20").Value = $true
$ws.Range("J21").Value = $true
$ws.Range("J26").Value = $true
$ws.Range("J27").Value = $true
$ws.Range("J32").Value = $true

$excel.Calculate()
